# feat: display note by id
#
# - NOTES sheet: remove the stray "c" row (old row 2) which shifts every
#   following note/tag pair up by one row so tags line up with the note
#   they actually belong to; append a new note "Clarity equals success".
# - ENGLISH sheet: append 10 new vocabulary entries (opulent, trait, tenet,
#   lure, proverb, insipid, startling, plateau, counterfeit, repercussion)
#   starting two rows below the last existing word (row 165 stays blank),
#   then make NOTES the active/selected sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ENGLISH")
$ws2 = $wb.Worksheets.Item("NOTES")

# --- NOTES: drop the stray "c" entry (row 2) ---------------------------
$ws2.Rows.Item(2).Delete()

# --- ENGLISH: append the newly studied words (rows 166-175) ------------
$ws1.Range("A166").Value = "opulent"
$ws1.Range("C166").Value = "luxurious;wealthy"
$ws1.Range("E166").Value = 0
$ws1.Range("F166").Value = "2021-11-21 15:33:36.099685"

$ws1.Range("A167").Value = "trait"
$ws1.Range("C167").Value = "characteristic"
$ws1.Range("E167").Value = 0
$ws1.Range("F167").Value = "2021-11-21 15:35:43.609799"

$ws1.Range("A168").Value = "tenet"
$ws1.Range("C168").Value = "principle"
$ws1.Range("E168").Value = 0
$ws1.Range("F168").Value = "2021-11-21 15:37:00.336208"

$ws1.Range("A169").Value = "lure"
$ws1.Range("C169").Value = "tempt"
$ws1.Range("E169").Value = 0
$ws1.Range("F169").Value = "2021-11-21 15:37:54.117744"

$ws1.Range("A170").Value = "proverb"
$ws1.Range("C170").Value = "saying"
$ws1.Range("E170").Value = 0
$ws1.Range("F170").Value = "2021-11-21 15:38:34.634497"

$ws1.Range("A171").Value = "insipid"
$ws1.Range("C171").Value = "tasteless;weak"
$ws1.Range("E171").Value = 0
$ws1.Range("F171").Value = "2021-11-21 15:40:27.314993"

$ws1.Range("A172").Value = "startling"
$ws1.Range("C172").Value = "remarkable;surprising;astonishing"
$ws1.Range("E172").Value = 0
$ws1.Range("F172").Value = "2021-11-21 15:41:25.484576"

$ws1.Range("A173").Value = "plateau"
$ws1.Range("B173").Value = "a state of little or no change following a period of activity or progress"
$ws1.Range("E173").Value = 0
$ws1.Range("F173").Value = "2021-11-21 15:42:35.633005"

$ws1.Range("A174").Value = "counterfeit"
$ws1.Range("C174").Value = "fake"
$ws1.Range("E174").Value = 0
$ws1.Range("F174").Value = "2021-11-21 15:44:40.04142"

$ws1.Range("A175").Value = "repercussion"
$ws1.Range("C175").Value = "consequence"
$ws1.Range("E175").Value = 0
$ws1.Range("F175").Value = "2021-11-21 15:45:55.726065"

# --- NOTES: append the new note at the bottom ---------------------------
$ws2.Range("A48").Value = "Clarity equals success"

# --- Selections / active sheet -----------------------------------------
$ws1.Range("H165").Select()
$ws2.Range("K29").Select()
$ws2.Activate()
